$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trialTypes")
$ws.Range("C2").Value = 13
$ws.Range("C3").Select()
